$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: swap the sample account/VM/script for the new ones ---
$ws.Range("A2").Value = "AZ-AS-SUB-EX-N-SEQ02125-CORE"
$ws.Range("B2").Value = "ubuntu18-script"
$ws.Range("C2").Value = "sudo apt-get update`nsudo apt-get upgrade -y`nsudo apt-get dist-upgrade -y`nsudo apt-get autoremove -y`nsudo apt-get install -y update-manager-core`nsudo sed -i 's/^Prompt=lts/Prompt=normal/' /etc/update-manager/release-upgrades`nexport DEBIAN_FRONTEND=noninteractive`nsudo -E do-release-upgrade -f DistUpgradeViewNonInteractive -m server --allow-third-party --allow-unauthenticated`ncat /etc/os-release"

# --- Rows 3-5: these were other sample rows; clear them out, leaving the formatted-but-empty cells ---
$ws.Range("A3:F5").ClearContents()

# --- Resize column C (now holding the long script) and refit all affected rows ---
$ws.Columns.Item(3).ColumnWidth = 62.6328125
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).EntireRow.AutoFit()
$ws.Rows.Item(4).EntireRow.AutoFit()
$ws.Rows.Item(5).EntireRow.AutoFit()

# --- Selection state ---
$ws.Range("D2").Select()
